# Update the "Förändrad" (Changed) date column C for all data rows (2-43)
# from 45839 (2025-07-01) to 45840 (2025-07-02).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45839) {
        $cell.Value = 45840
    }
}
